$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60 (shifts existing rows 60.. down to 61..)
$ws.Rows.Item(60).Insert()

# Fill the new row 60 with data (same pattern as neighboring rows)
$ws.Cells.Item(60, 1).Value = 6
$ws.Cells.Item(60, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(60, 3).Value = "Metropolitana"
$ws.Cells.Item(60, 4).Value = 44494
$ws.Cells.Item(60, 5).Value = 13
$ws.Cells.Item(60, 6).Value = 100112026
$ws.Cells.Item(60, 7).Value = "Haba"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Primera"
$ws.Cells.Item(60, 10).Value = 980
$ws.Cells.Item(60, 11).Value = 5000
$ws.Cells.Item(60, 12).Value = 6000
$ws.Cells.Item(60, 13).Value = 5571
$ws.Cells.Item(60, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(60, 15).Value = "Región Metropolitana"
$ws.Cells.Item(60, 16).Value = 223
$ws.Cells.Item(60, 17).Value = 25
$ws.Cells.Item(60, 18).Value = "Hortaliza"

Write-Output "Done"
